$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data refresh.
# Some Price values are plain decimal numbers that Excel would otherwise auto-convert
# to a numeric type; force those to remain text (matching the original inlineStr cells)
# and then restore the default "Normal" style so no stray number-format style is left behind.
$ws.Range("D2").Value = "66.742.10"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.523.84"
$ws.Range("E3").Value = "  -3.89%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "585.06"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "171.56"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "2.521.78"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("E13").Value = "  -1.75%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.78"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.16%  "
$ws.Range("D15").Value = "2.983.14"
$ws.Range("E15").Value = "  -4.31%  "
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "66.590.53"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "2.522.01"
$ws.Range("E18").Value = "  -4.44%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.86"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.34"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -5.60%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "348.05"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.42%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.96"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.97%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "2.634.48"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "0.0₃0982"
$ws.Range("E30").Value = "  -2.58%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "528.06"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("E33").Value = "  -2.48%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  -2.63%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "157.09"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("E44").Value = "  -0.03%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.52"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("E46").Value = "  -1.82%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "149.51"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("E48").Value = "  -3.62%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "0.0₆0270"
$ws.Range("E51").Value = "  -10.56%  "
